# Apply the "Sprint Backlog" sheet edit:
#  - delete the 3 blank rows (4:6) that sat above the task table, shifting
#    everything below up by three rows
#  - the task row that used to be a spacer (now row 11) loses its stray
#    "x" flag and leftover 0 in column D
#  - the "Meteo" task row (now row 12) gains the "x" flag that used to sit
#    on that spacer row
#  - tidy up the view state (no frozen/scrolled topLeftCell, new selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")

$ws.Rows("4:6").Delete()

# Row 11 (used to be row 14): clear the leftover value/flag cells so only
# the empty-but-styled E:H cells remain.
$ws.Range("D11").Clear()
$ws.Range("J11").Clear()
$ws.Range("E11:H11").ClearContents()

# Row 12 (used to be row 15, "Meteo (Jour, nuit, pluie)"): flag it as
# "En attente" the same way its neighbours are flagged.
$ws.Range("J12").Value = "x"

$ws.Activate()
$ws.Range("P11").Select()
